$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.23878994128646
$ws.Range("D2").Value = 2.99935422207211
$ws.Range("E2").Value = 21.56443023777952
$ws.Range("F2").Value = 18.19196770799345
$ws.Range("G2").Value = 19.46581094307508
$ws.Range("H2").Value = 11.18709037342399
$ws.Range("L2").Value = 8.479934202105618
$ws.Range("N2").Value = 18.36521348911732
$ws.Range("O2").Value = 15.84820985431647
$ws.Range("B3").Value = 16.01044332756856
$ws.Range("D3").Value = 2.954533584693999
$ws.Range("E3").Value = 21.58934862633505
$ws.Range("F3").Value = 17.98723077660098
$ws.Range("G3").Value = 18.98239725730938
$ws.Range("H3").Value = 11.18875262248342
$ws.Range("L3").Value = 8.317447209622369
$ws.Range("N3").Value = 18.35927307079072
$ws.Range("O3").Value = 15.76984899705695
$ws.Range("B4").Value = 15.8712669572278
$ws.Range("D4").Value = 2.93556071399322
$ws.Range("E4").Value = 21.60792589587436
$ws.Range("F4").Value = 17.86642068448433
$ws.Range("G4").Value = 18.68709121071567
$ws.Range("H4").Value = 11.19219014927639
$ws.Range("L4").Value = 8.216904228700527
$ws.Range("N4").Value = 18.35803865638474
$ws.Range("O4").Value = 15.72624861026139
$ws.Range("B5").Value = 15.81487664396024
$ws.Range("D5").Value = 2.927910546689169
$ws.Range("E5").Value = 21.61632420022625
$ws.Range("F5").Value = 17.81848330814386
$ws.Range("G5").Value = 18.56735113869991
$ws.Range("H5").Value = 11.19419908588276
$ws.Range("L5").Value = 8.175787946256657
$ws.Range("N5").Value = 18.35814569460835
$ws.Range("O5").Value = 15.7096315651265
$ws.Range("B6").Value = 15.8055346509957
$ws.Range("D6").Value = 2.926645388520115
$ws.Range("E6").Value = 21.61776884521066
$ws.Range("F6").Value = 17.81060322669881
$ws.Range("G6").Value = 18.54751111514608
$ws.Range("H6").Value = 11.19456940175338
$ws.Range("L6").Value = 8.168953419035232
$ws.Range("N6").Value = 18.35820040771949
$ws.Range("O6").Value = 15.70694221959301
$ws.Range("B7").Value = 15.87050505357
$ws.Range("D7").Value = 2.935457201112191
$ws.Range("E7").Value = 21.60803580119202
$ws.Range("F7").Value = 17.86576886667629
$ws.Range("G7").Value = 18.68547363293857
$ws.Range("H7").Value = 11.19221478001174
$ws.Range("L7").Value = 8.216350235951948
$ws.Range("N7").Value = 18.35803762556621
$ws.Range("O7").Value = 15.72601982968991
$ws.Range("B8").Value = 16.15988220435825
$ws.Range("D8").Value = 2.980964191220068
$ws.Range("E8").Value = 21.57234423141502
$ws.Range("F8").Value = 18.12039553093512
$ws.Range("G8").Value = 19.29894588567738
$ws.Range("H8").Value = 11.18716219812285
$ws.Range("L8").Value = 8.424099612601481
$ws.Range("N8").Value = 18.36266643942652
$ws.Range("O8").Value = 15.82026412085769
$ws.Range("B9").Value = 16.73260862630778
$ws.Range("D9").Value = 3.109800199599588
$ws.Range("E9").Value = 21.52818387215578
$ws.Range("F9").Value = 18.65574193768736
$ws.Range("G9").Value = 20.50466100833032
$ws.Range("H9").Value = 11.19640309373429
$ws.Range("L9").Value = 8.823191976602111
$ws.Range("N9").Value = 18.39073899288409
$ws.Range("O9").Value = 16.04012122102209
$ws.Range("B10").Value = 17.15274225677634
$ws.Range("D10").Value = 3.199047305325468
$ws.Range("E10").Value = 21.51124607917591
$ws.Range("F10").Value = 19.06706029028471
$ws.Range("G10").Value = 21.38032324036154
$ws.Range("H10").Value = 11.21481089851096
$ws.Range("L10").Value = 9.108591868378269
$ws.Range("N10").Value = 18.42273287580949
$ws.Range("O10").Value = 16.22186710312983
$ws.Range("B11").Value = 17.34293633600128
$ws.Range("D11").Value = 3.238381378616649
$ws.Range("E11").Value = 21.50685352695269
$ws.Range("F11").Value = 19.25722048263967
$ws.Range("G11").Value = 21.7741915611853
$ws.Range("H11").Value = 11.22568948887934
$ws.Range("L11").Value = 9.236191297578946
$ws.Range("N11").Value = 18.43970749995704
$ws.Range("O11").Value = 16.30866086781238
$ws.Range("B12").Value = 17.41476050537089
$ws.Range("D12").Value = 3.253087158623459
$ws.Range("E12").Value = 21.5056617095149
$ws.Range("F12").Value = 19.32959466501976
$ws.Range("G12").Value = 21.92252254875454
$ws.Range("H12").Value = 11.23016696965934
$ws.Range("L12").Value = 9.28414733489292
$ws.Range("N12").Value = 18.44647890846276
$ws.Range("O12").Value = 16.34209402797982
$ws.Range("B13").Value = 17.39930171096861
$ws.Range("D13").Value = 3.24992853062147
$ws.Range("E13").Value = 21.50589749454998
$ws.Range("F13").Value = 19.31399254326664
$ws.Range("G13").Value = 21.89061586075724
$ws.Range("H13").Value = 11.22918678061873
$ws.Range("L13").Value = 9.273835940487533
$ws.Range("N13").Value = 18.44500536268584
$ws.Range("O13").Value = 16.33486883675449
$ws.Range("B14").Value = 17.34884969003271
$ws.Range("D14").Value = 3.23959506350119
$ws.Range("E14").Value = 21.50674605292663
$ws.Range("F14").Value = 19.26316786361371
$ws.Range("G14").Value = 21.78641219803724
$ws.Range("H14").Value = 11.22605069291872
$ws.Range("L14").Value = 9.24014418024082
$ws.Range("N14").Value = 18.44025773070223
$ws.Range("O14").Value = 16.31140023710947
$ws.Range("B15").Value = 17.31791866889306
$ws.Range("D15").Value = 3.233240671471887
$ws.Range("E15").Value = 21.50732708381517
$ws.Range("F15").Value = 19.23208158256845
$ws.Range("G15").Value = 21.72247289848703
$ws.Range("H15").Value = 11.224176299348
$ws.Range("L15").Value = 9.219458504408802
$ws.Range("N15").Value = 18.43739426213543
$ws.Range("O15").Value = 16.29709800055932
$ws.Range("B16").Value = 17.14028837792301
$ws.Range("D16").Value = 3.196450644256937
$ws.Range("E16").Value = 21.51159939349722
$ws.Range("F16").Value = 19.05468791488822
$ws.Range("G16").Value = 21.35447772308844
$ws.Range("H16").Value = 11.21415018022683
$ws.Range("L16").Value = 9.100204384384138
$ws.Range("N16").Value = 18.42167186358144
$ws.Range("O16").Value = 16.21627559233516
$ws.Range("B17").Value = 17.03103370870868
$ws.Range("D17").Value = 3.173551847700016
$ws.Range("E17").Value = 21.51506537810313
$ws.Range("F17").Value = 18.94659103119165
$ws.Range("G17").Value = 21.12744989657562
$ws.Range("H17").Value = 11.20863956582097
$ws.Range("L17").Value = 9.026443634731958
$ws.Range("N17").Value = 18.41264342665955
$ws.Range("O17").Value = 16.16773011844293
$ws.Range("B18").Value = 16.968109834918
$ws.Range("D18").Value = 3.160262414274136
$ws.Range("E18").Value = 21.51737100900228
$ws.Range("F18").Value = 18.88470800343912
$ws.Range("G18").Value = 20.99645794570753
$ws.Range("H18").Value = 11.2057059211982
$ws.Range("L18").Value = 8.983810840375444
$ws.Range("N18").Value = 18.40767861845782
$ws.Range("O18").Value = 16.14019714539617
$ws.Range("B19").Value = 16.94679256790281
$ws.Range("D19").Value = 3.155742691992633
$ws.Range("E19").Value = 21.5182054008459
$ws.Range("F19").Value = 18.86380778443303
$ws.Range("G19").Value = 20.9520411817145
$ws.Range("H19").Value = 11.20475321939671
$ws.Range("L19").Value = 8.96934179220667
$ws.Range("N19").Value = 18.40603693383204
$ws.Range("O19").Value = 16.13094254611557
$ws.Range("B20").Value = 17.04267317146389
$ws.Range("D20").Value = 3.176001797975777
$ws.Range("E20").Value = 21.514664152375
$ws.Range("F20").Value = 18.9580685021587
$ws.Range("G20").Value = 21.15166129876872
$ws.Range("H20").Value = 11.20920177859617
$ws.Range("L20").Value = 9.034317389303126
$ws.Range("N20").Value = 18.4135809431412
$ws.Range("O20").Value = 16.17285777921899
$ws.Range("B21").Value = 17.36367455730229
$ws.Range("D21").Value = 3.242635443308585
$ws.Range("E21").Value = 21.50648405209547
$ws.Range("F21").Value = 19.27808699406774
$ws.Range("G21").Value = 21.81704287982143
$ws.Range("H21").Value = 11.22696214115371
$ws.Range("L21").Value = 9.250050446268803
$ws.Range("N21").Value = 18.44164294056792
$ws.Range("O21").Value = 16.31827837676248
$ws.Range("B22").Value = 17.57228372576059
$ws.Range("D22").Value = 3.285078650844219
$ws.Range("E22").Value = 21.50388408412795
$ws.Range("F22").Value = 19.48932917025437
$ws.Range("G22").Value = 22.24706610320398
$ws.Range("H22").Value = 11.24065497809721
$ws.Range("L22").Value = 9.388910219803162
$ws.Range("N22").Value = 18.46198292358547
$ws.Range("O22").Value = 16.4166069343573
$ws.Range("B23").Value = 17.4610742792575
$ws.Range("D23").Value = 3.262529372405395
$ws.Range("E23").Value = 21.50502208251076
$ws.Range("F23").Value = 19.37641812480426
$ws.Range("G23").Value = 22.01805206390177
$ws.Range("H23").Value = 11.23315684865763
$ws.Range("L23").Value = 9.315006764420957
$ws.Range("N23").Value = 18.45094566639442
$ws.Range("O23").Value = 16.36383509106927
$ws.Range("B24").Value = 17.03741131392811
$ws.Range("D24").Value = 3.174894563031085
$ws.Range("E24").Value = 21.51484457124293
$ws.Range("F24").Value = 18.95287871169229
$ws.Range("G24").Value = 21.14071677625557
$ws.Range("H24").Value = 11.20894687163106
$ws.Range("L24").Value = 9.030758367911595
$ws.Range("N24").Value = 18.41315638837099
$ws.Range("O24").Value = 16.17053838915574
$ws.Range("B25").Value = 16.57750978705801
$ws.Range("D25").Value = 3.075863529644245
$ws.Range("E25").Value = 21.53738925319526
$ws.Range("F25").Value = 18.50746315244938
$ws.Range("G25").Value = 20.17945530179186
$ws.Range("H25").Value = 11.19185755356283
$ws.Range("L25").Value = 8.716407985545459
$ws.Range("N25").Value = 18.38113203130776
$ws.Range("O25").Value = 15.97700430481967
